$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Espèces")

# Insert a new column before column F (shifts F:K to G:L)
$ws.Range("F1").EntireColumn.Insert()

# Header for the new "Renomée" column
$ws.Range("F1").Value = "Renomée"

# New Renomée values on specific rows
$ws.Range("F4").Value = "Sagesse"
$ws.Range("F10").Value = "Ruse"
$ws.Range("F13").Value = "Gloire"
